# Weekly price update: a new daily reading for "Macroferia Regional de Talca -
# Zapallo italiano" is inserted above the existing row 258, pushing every
# following record (258-275) down by one row (now 259-276) and extending the
# sheet's used range to R276.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 258; rows 258:275 shift down to 259:276.
$ws.Rows("258:258").Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(258, 1).Value2 = 5
$ws.Cells.Item(258, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(258, 3).Value2 = "Maule"
$ws.Cells.Item(258, 4).Value2 = 44585
$ws.Cells.Item(258, 5).Value2 = 7
$ws.Cells.Item(258, 6).Value2 = 100112032
$ws.Cells.Item(258, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(258, 8).Value2 = "Sin especificar"
$ws.Cells.Item(258, 9).Value2 = "Primera"
$ws.Cells.Item(258, 10).Value2 = 300
$ws.Cells.Item(258, 11).Value2 = 8000
$ws.Cells.Item(258, 12).Value2 = 8000
$ws.Cells.Item(258, 13).Value2 = 8000
$ws.Cells.Item(258, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(258, 15).Value2 = "Región del Maule"
$ws.Cells.Item(258, 16).Value2 = 133
$ws.Cells.Item(258, 17).Value2 = 60
$ws.Cells.Item(258, 18).Value2 = "Hortaliza"
